# FilterData sheet: update the second test-data row (TestData\DVV -> DEV
# rerun). FilterValue1 (F2) and Operator2 (H2) are the two data cells the
# commit actually changes; selection is returned to A1 afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FilterValue1: 233 -> 23
$ws.Range("F2").Value = 23

# Operator2: "Equal" -> "Does Not Contains"
$ws.Range("H2").Value = "Does Not Contains"

# Return to the default selection (A1) instead of leaving A2 selected.
$ws.Range("A1").Select()
